$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rundata")

$rows = @(
    @("QA","Sub-D","deluxe25offp-redes-spring","Kit",30,"BHK738","https://example.com","Same","Visa","Chrome"),
    @("QA","Sub-D","deluxe25offp-redes-spring","Kit",90,"BHK73A","https://example.com","Same","Paypal","Chrome")
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
Write-Output $ws.Cells.Item(2,1).Value()
Write-Output $ws.Cells.Item(2,5).Value()
Write-Output $ws.Cells.Item(3,6).Value()
